# Updated mint- and transfer gas consumption data of the ERC721F contract
# (columns I:L, "ERC721F" block) on the three data sheets, then restore
# the sheet/selection view state to match the author's final session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Optimizer Disabled" -> ERC721F block (I3:L10)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Optimizer Disabled")
$ws1.Range("J4").Value = 74559
$ws1.Range("K4").Value = 305347
$ws1.Range("L4").Value = 2613690
$ws1.Range("J5").Value = 60149
$ws1.Range("K5").Value = 64949
$ws1.Range("J6").Value = 60005
$ws1.Range("K6").Value = 65167
$ws1.Range("K7").Value = 207212
$ws1.Range("L7").Value = 212012
$ws1.Range("K8").Value = 208944
$ws1.Range("L8").Value = 214106
$ws1.Range("L9").Value = 1258333
$ws1.Range("L10").Value = 1271818

# ---------------------------------------------------------------------
# Sheet "Runs 200 - Optimizer Enabled" -> ERC721F block (I3:L10)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Runs 200 - Optimizer Enabled")
$ws2.Range("J4").Value = 73891
$ws2.Range("K4").Value = 300206
$ws2.Range("L4").Value = 2563819
$ws2.Range("J5").Value = 59329
$ws2.Range("K5").Value = 64129
$ws2.Range("J6").Value = 59135
$ws2.Range("K6").Value = 64217
$ws2.Range("K7").Value = 196776
$ws2.Range("L7").Value = 201576
$ws2.Range("K8").Value = 197018
$ws2.Range("L8").Value = 202100
$ws2.Range("L9").Value = 1129757
$ws2.Range("L10").Value = 1129392

# ---------------------------------------------------------------------
# Sheet "Runs 1000 - Optimizer Enabled" -> ERC721F block (I3:L10)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Runs 1000 - Optimizer Enabled")
$ws3.Range("J4").Value = 73879
$ws3.Range("K4").Value = 300086
$ws3.Range("L4").Value = 2562619
$ws3.Range("J5").Value = 59305
$ws3.Range("K5").Value = 64105
$ws3.Range("J6").Value = 59111
$ws3.Range("K6").Value = 64193
$ws3.Range("K7").Value = 196536
$ws3.Range("L7").Value = 201336
$ws3.Range("K8").Value = 196778
$ws3.Range("L8").Value = 201860
$ws3.Range("L9").Value = 1128557
$ws3.Range("L10").Value = 1128192

$ws4 = $wb.Worksheets.Item("Graphs")

# ---------------------------------------------------------------------
# Restore view / selection state for each sheet, matching the author's
# last-saved session (tab switched from "Graphs" back to "Optimizer
# Disabled", scroll position and selected cell changed on every sheet).
# ---------------------------------------------------------------------
$ws2.Select()
$ws2.Range("H13").Select()

$ws3.Select()
$ws3.Range("D1").Select()
$ws3.Range("M9").Select()

$ws4.Select()
$ws4.Range("AK11").Select()

$ws1.Select()
$ws1.Range("N14").Select()
